$d = $word.ActiveDocument

# Locate the paragraph that holds the "Ver no Jupiter Salvar em pdf Salvar em docx"
# text. In the surrounding layout it is preceded by a blank paragraph and
# followed by a blank paragraph and then a page-break paragraph; all four of
# those paragraphs are removed, leaving the blank paragraph and page-break
# paragraph that follow them untouched.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $startPar = $d.Paragraphs.Item($target - 1)
    $afterPar = $d.Paragraphs.Item($target + 3)
    $range = $d.Range($startPar.Range.Start, $afterPar.Range.Start)
    $range.Delete()
}
